$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

$ws.Range("A17").Value = "30/01/2018"
$ws.Range("B17").Value = "2222"
$ws.Range("C17").Value = "Гострини на розрізі контакту"
$ws.Range("D17").Value = "2"

$ws.Range("A18").Value = "31/01/2018"
$ws.Range("B18").Value = "3012"
$ws.Range("C18").Value = "Гострини на розрізі контакту"
$ws.Range("D18").Value = "1"

$ws.Range("A19").Value = "31/01/2018"
$ws.Range("B19").Value = "3012"
$ws.Range("C19").Value = "Не вірна довжина проводу"
$ws.Range("D19").Value = "1"

$ws.Range("A20").Value = "31/01/2018"
$ws.Range("B20").Value = "3012"
$ws.Range("C20").Value = "Не вірна довжина проводу"
$ws.Range("D20").Value = "123654789"

$ws.Range("A21").Value = "31/01/2018"
$ws.Range("B21").Value = "3012"
$ws.Range("C21").Value = "інше"
$ws.Range("D21").Value = "123456"

$ws.Range("A22").Value = "**"
